# Apply weekly data refresh: the rows of daily price records are
# re-associated with a different rotation of dates/prices (same data set,
# shuffled across the existing rows). Columns D (Fecha), L (Calidad),
# M (Volumen), N (Precio minimo), O (Precio maximo), P (Precio promedio
# ponderado), Q (Unidad de comercializacion), S (Precio $/Kg) and
# T (Kg / unidad) move between rows 2-20 according to the mapping below;
# all other columns are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns (by index) that participate in the row permutation.
$cols = @(4, 12, 13, 14, 15, 16, 17, 19, 20)

# Snapshot the "before" values for every relevant cell so the permutation
# (which includes multi-row cycles) can be applied without clobbering
# data that is still needed as a source.
$snapshot = @{}
foreach ($r in 2..20) {
    foreach ($c in $cols) {
        $snapshot["$r,$c"] = $ws.Cells.Item($r, $c).Value2
    }
}

# Maps destination row -> source row (source row's "before" values become
# the destination row's "after" values).
$rowMap = @{
    2  = 19
    3  = 4
    4  = 2
    5  = 13
    6  = 7
    7  = 11
    8  = 5
    9  = 6
    10 = 15
    11 = 14
    12 = 20
    13 = 10
    14 = 18
    15 = 16
    16 = 17
    17 = 8
    18 = 9
    19 = 3
    20 = 12
}

foreach ($destRow in 2..20) {
    $srcRow = $rowMap[$destRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($destRow, $c).Value = $snapshot["$srcRow,$c"]
    }
}
